# "Fin sprite punto de mira"
# Adds the new "punto mira" sheet (the crosshair sprite sheet) at the end of
# the workbook, fills in its sprite byte-grid + BIN2HEX formulas, makes it
# the active sheet/tab, and tidies up the selections left on the two sheets
# ("jefezombid" / "jefezombii") that were being worked on before this sheet
# was finished.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Fix up the selection left on "jefezombii" (the previously-active sheet).
#    "jefezombid" keeps its existing selection (O17:P17 / activeCell P17) -
#    it just stops being the selected tab once we activate the new sheet
#    below, so it is left untouched here.
# ---------------------------------------------------------------------------
$jefezombii = $wb.Worksheets.Item("jefezombii")
$jefezombii.Range("C1").Select()

# ---------------------------------------------------------------------------
# 2. Add the new sheet at the end of the workbook.
# ---------------------------------------------------------------------------
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "punto mira"

# ---------------------------------------------------------------------------
# 3. Sprite data: pairs of binary-string cells (A/B, E/F, I/J, M/N) each
#    followed by a BIN2HEX-based hex formula (C/D, G/H, K/L, O/P).
# ---------------------------------------------------------------------------
$data = @{}
$data["A"] = @("n00000000","n00000000","n00000000","n00000001","n00000010","n00000100","n00000001","n00001010","n00001010","n00000001","n00000100","n00000010","n00000001","n00000000","n00000000","n00000000")
$data["B"] = @("n00000000","n00000000","n00000000","n10000000","n01000000","n00100000","n10000000","n01010000","n01010000","n10000000","n00100000","n01000000","n10000000","n00000000","n00000000","n00000000")
$data["E"] = @("n00000000","n00000000","n00000000","n00000010","n00000101","n00001011","n00001110","n00000101","n00000101","n00001110","n00001011","n00000101","n00000010","n00000000","n00000000","n00000000")
$data["F"] = @("n00000000","n00000000","n00000000","n01000000","n10100000","n11010000","n01110000","n10100000","n10100000","n01110000","n11010000","n10100000","n01000000","n00000000","n00000000","n00000000")
$data["I"] = @("n00000000","n00000000","n00000000","n00000010","n00000100","n00001100","n00010010","n00000001","n00000001","n00010010","n00001100","n00000100","n00000010","n00000000","n00000000","n00000000")
$data["J"] = @("n00000000","n00000000","n00000000","n01000000","n00100000","n00110000","n01001000","n10000000","n10000000","n01001000","n00110000","n00100000","n01000000","n00000000","n00000000","n00000000")
$data["M"] = @("n00000000","n00000000","n00000001","n00000001","n00000001","n00000011","n00000101","n00111111","n00111111","n00000101","n00000011","n00000001","n00000001","n00000001","n00000000","n00000000")
$data["N"] = @("n00000000","n00000000","n10000000","n10000000","n10000000","n11000000","n10100000","n11111100","n11111100","n10100000","n11000000","n10000000","n10000000","n10000000","n00000000","n00000000")

foreach ($col in @("A","B","E","F","I","J","M","N")) {
    $vals = $data[$col]
    for ($i = 0; $i -lt 16; $i++) {
        $ws.Range($col + ($i + 1)).Value = $vals[$i]
    }
}

# Hex-formula columns, keyed by the source binary column they read from.
$formulaCols = [ordered]@{
    "C" = "A"
    "D" = "B"
    "G" = "E"
    "H" = "F"
    "K" = "I"
    "L" = "J"
    "O" = "M"
    "P" = "N"
}

foreach ($dest in $formulaCols.Keys) {
    $src = $formulaCols[$dest]
    $ws.Range($dest + "1").Formula = '=CONCATENATE("#",BIN2HEX(REPLACE(' + $src + '1,1,1,""),2))'
    $formula = '=CONCATENATE("#",BIN2HEX(REPLACE(' + $src + '2,1,1,""),2))'
    $ws.Range($dest + "2:" + $dest + "16").Formula = $formula
}

# ---------------------------------------------------------------------------
# 4. Selection + activation: "punto mira" becomes the active/selected tab.
# ---------------------------------------------------------------------------
$ws.Range("O1:P16").Select()
$ws.Activate()
